$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the recalculated results in row 2 (columns O through AI) with the
# refreshed values produced by the workbook's heat-exchanger network solver.
$ws.Range("O2").Value  = 0.9999999998517366
$ws.Range("P2").Value  = 0.9999999998531801
$ws.Range("Q2").Value  = 20002999.99999954
$ws.Range("R2").Value  = 20001999.99999954
$ws.Range("S2").Value  = 20000999.99999977
$ws.Range("T2").Value  = 20000000
$ws.Range("U2").Value  = 789235.1064058846
$ws.Range("V2").Value  = 915985.2593332215
$ws.Range("W2").Value  = 956727.0286293968
$ws.Range("X2").Value  = 969979.3779034525
$ws.Range("Y2").Value  = 5
$ws.Range("Z2").Value  = 2.500000000432172
$ws.Range("AA2").Value = 2.500000000432172
$ws.Range("AB2").Value = 132575.0000016871
$ws.Range("AC2").Value = 107575.0000016871
$ws.Range("AD2").Value = 101325
$ws.Range("AE2").Value = 101325
$ws.Range("AF2").Value = 918999.5861496417
$ws.Range("AG2").Value = 916349.1160479402
$ws.Range("AH2").Value = 900052.4083350284
$ws.Range("AI2").Value = 865649.0549721529

# The stray formatted-but-empty cell below the table (row 5) is no longer
# part of the sheet, so remove the whole row and let the used range shrink
# back down to the real data (A1:AI2).
$ws.Rows.Item(5).Delete()

# Reset the view back to the top-left of the sheet instead of the scrolled,
# mid-table selection that was saved previously.
$ws.Range("A1").Select()
